$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$s.Shapes.Item("Textfeld 1").Delete()
